$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- describe() table header row (row 2) ---
$ws.Range("I2").Value = "count"
$ws.Range("J2").Value = "mean"
$ws.Range("K2").Value = "std"
$ws.Range("L2").Value = "min"
$ws.Range("M2").Value = 0.25
$ws.Range("N2").Value = 0.5
$ws.Range("O2").Value = 0.75
$ws.Range("P2").Value = "max"
$ws.Range("M2:O2").NumberFormat = "0%"

# --- describe() table body (rows 3-7), one row per group ---
$ws.Range("H3").Value = "A型号"
$ws.Range("I3").Value = 10
$ws.Range("J3").Value = 274.39999999999998
$ws.Range("K3").Value = 6.6030296076876711
$ws.Range("L3").Value = 265
$ws.Range("M3").Value = 268.75
$ws.Range("N3").Value = 273.5
$ws.Range("O3").Value = 279.75
$ws.Range("P3").Value = 285

$ws.Range("H4").Value = "B型号"
$ws.Range("I4").Value = 10
$ws.Range("J4").Value = 282.89999999999998
$ws.Range("K4").Value = 27.79868102386634
$ws.Range("L4").Value = 263
$ws.Range("M4").Value = 268.75
$ws.Range("N4").Value = 277
$ws.Range("O4").Value = 283.25
$ws.Range("P4").Value = 359

$ws.Range("H5").Value = "C型号"
$ws.Range("I5").Value = 10
$ws.Range("J5").Value = 262.7
$ws.Range("K5").Value = 5.6184220797895446
$ws.Range("L5").Value = 254
$ws.Range("M5").Value = 259
$ws.Range("N5").Value = 261.5
$ws.Range("O5").Value = 267.75
$ws.Range("P5").Value = 271

$ws.Range("H6").Value = "D型号"
$ws.Range("I6").Value = 10
$ws.Range("J6").Value = 268
$ws.Range("K6").Value = 10.934146311237816
$ws.Range("L6").Value = 258
$ws.Range("M6").Value = 261
$ws.Range("N6").Value = 262.5
$ws.Range("O6").Value = 270.75
$ws.Range("P6").Value = 290

$ws.Range("H7").Value = "E型号"
$ws.Range("I7").Value = 10
$ws.Range("J7").Value = 285
$ws.Range("K7").Value = 12.092238098144703
$ws.Range("L7").Value = 259
$ws.Range("M7").Value = 284.75
$ws.Range("N7").Value = 288.5
$ws.Range("O7").Value = 293.25
$ws.Range("P7").Value = 295

# --- ANOVA table ---
$ws.Range("H14").Value = "方差分析"

$ws.Range("I15").Value = "sum_sq"
$ws.Range("J15").Value = "df"
$ws.Range("K15").Value = "F"
$ws.Range("L15").Value = "PR(>F)"

$ws.Range("H16").Value = "Intercept"
$ws.Range("I16").Value = 752953.60000000114
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 3380.3811082068014
$ws.Range("L16").Value = 5.5207908951615785 * [Math]::Pow(10, -44)

$ws.Range("H17").Value = "C(Threat)"
$ws.Range("I17").Value = 3622.599999999989
$ws.Range("J17").Value = 4
$ws.Range("K17").Value = 4.0659107688009932
$ws.Range("L17").Value = 6.7397238407765434 * [Math]::Pow(10, -3)

$ws.Range("H18").Value = "Residual"
$ws.Range("I18").Value = 10023.4
$ws.Range("J18").Value = 45
